$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

$ws.Cells.Item($row, 1).Value = 42632.883483796293
$ws.Cells.Item($row, 2).Value = 14
$ws.Cells.Item($row, 3).Value = "Buy"
$ws.Cells.Item($row, 4).Value = 26
$ws.Cells.Item($row, 5).Value = 10694
$ws.Cells.Item($row, 6).Value = 1727
$ws.Cells.Item($row, 7).Value = 63
$ws.Cells.Item($row, 8).Value = 35
$ws.Cells.Item($row, 9).Value = 83
$ws.Cells.Item($row, 10).Value = 16
$ws.Cells.Item($row, 11).Value = 10004
$ws.Cells.Item($row, 12).Value = 245
$ws.Cells.Item($row, 13).Value = 138
$ws.Cells.Item($row, 14).Value = 30
$ws.Cells.Item($row, 15).Value = 6
$ws.Cells.Item($row, 16).Value = "Noun"
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = 0.87
$ws.Cells.Item($row, 19).Value = 0.0351
$ws.Cells.Item($row, 20).Value = -2.08
$ws.Cells.Item($row, 21).Value = 15.16
$ws.Cells.Item($row, 22).Value = "N/A"
$ws.Cells.Item($row, 23).Value = 0

$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("S3").NumberFormat = "0.00%"
